$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Resolving-Mac" sending-cluster rows (rows 14-17); the
# remaining data will be rewritten in place with the refreshed TPM values.
$ws.Rows.Item(14).Resize(4).Delete()

# Refresh every data cell (rows 2-13) with the recalculated TPM-based values.
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf13"
$ws.Cells.Item(2, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.342908333333333
$ws.Cells.Item(2, 8).Value = 4.028725
$ws.Cells.Item(2, 9).Value = 0.2879023314891748
$ws.Cells.Item(2, 10).Value = 0.2879023314891748
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.2347256666666667
$ws.Cells.Item(2, 14).Value = 0.7041770000000001
$ws.Cells.Item(2, 15).Value = 0.02963901389354162
$ws.Cells.Item(2, 16).Value = 0.02963901389354162
$ws.Cells.Item(2, 17).Value = 0.3152150538138889
$ws.Cells.Item(2, 18).Value = 2.836935484325
$ws.Cells.Item(2, 19).Value = 0.008533141202990677
$ws.Cells.Item(2, 20).Value = 0.008533141202990677
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf13"
$ws.Cells.Item(3, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.342908333333333
$ws.Cells.Item(3, 8).Value = 4.028725
$ws.Cells.Item(3, 9).Value = 0.2879023314891748
$ws.Cells.Item(3, 10).Value = 0.2879023314891748
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.035534
$ws.Cells.Item(3, 14).Value = 0.106602
$ws.Cells.Item(3, 15).Value = 0.004486909057068498
$ws.Cells.Item(3, 16).Value = 0.004486909057068498
$ws.Cells.Item(3, 17).Value = 0.04771890471666666
$ws.Cells.Item(3, 18).Value = 0.42947014245
$ws.Cells.Item(3, 19).Value = 0.001291791578709915
$ws.Cells.Item(3, 20).Value = 0.001291791578709915
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf13"
$ws.Cells.Item(4, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.342908333333333
$ws.Cells.Item(4, 8).Value = 4.028725
$ws.Cells.Item(4, 9).Value = 0.2879023314891748
$ws.Cells.Item(4, 10).Value = 0.2879023314891748
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.649223333333333
$ws.Cells.Item(4, 14).Value = 22.94767
$ws.Cells.Item(4, 15).Value = 0.9658740770493899
$ws.Cells.Item(4, 16).Value = 0.9658740770493899
$ws.Cells.Item(4, 17).Value = 10.27220575786111
$ws.Cells.Item(4, 18).Value = 92.44985182074998
$ws.Cells.Item(4, 19).Value = 0.2780773987074742
$ws.Cells.Item(4, 20).Value = 0.2780773987074742
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnfsf13"
$ws.Cells.Item(5, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.4963216666666666
$ws.Cells.Item(5, 8).Value = 1.488965
$ws.Cells.Item(5, 9).Value = 0.1064050028249084
$ws.Cells.Item(5, 10).Value = 0.1064050028249084
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.2347256666666667
$ws.Cells.Item(5, 14).Value = 0.7041770000000001
$ws.Cells.Item(5, 15).Value = 0.02963901389354162
$ws.Cells.Item(5, 16).Value = 0.02963901389354162
$ws.Cells.Item(5, 17).Value = 0.1164994340894444
$ws.Cells.Item(5, 18).Value = 1.048494906805
$ws.Cells.Item(5, 19).Value = 0.003153739357069796
$ws.Cells.Item(5, 20).Value = 0.003153739357069796
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf13"
$ws.Cells.Item(6, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.4963216666666666
$ws.Cells.Item(6, 8).Value = 1.488965
$ws.Cells.Item(6, 9).Value = 0.1064050028249084
$ws.Cells.Item(6, 10).Value = 0.1064050028249084
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.035534
$ws.Cells.Item(6, 14).Value = 0.106602
$ws.Cells.Item(6, 15).Value = 0.004486909057068498
$ws.Cells.Item(6, 16).Value = 0.004486909057068498
$ws.Cells.Item(6, 17).Value = 0.01763629410333333
$ws.Cells.Item(6, 18).Value = 0.15872664693
$ws.Cells.Item(6, 19).Value = 0.0004774295708924807
$ws.Cells.Item(6, 20).Value = 0.0004774295708924807
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf13"
$ws.Cells.Item(7, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.4963216666666666
$ws.Cells.Item(7, 8).Value = 1.488965
$ws.Cells.Item(7, 9).Value = 0.1064050028249084
$ws.Cells.Item(7, 10).Value = 0.1064050028249084
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.649223333333333
$ws.Cells.Item(7, 14).Value = 22.94767
$ws.Cells.Item(7, 15).Value = 0.9658740770493899
$ws.Cells.Item(7, 16).Value = 0.9658740770493899
$ws.Cells.Item(7, 17).Value = 3.796475273505555
$ws.Cells.Item(7, 18).Value = 34.16827746155
$ws.Cells.Item(7, 19).Value = 0.1027738338969461
$ws.Cells.Item(7, 20).Value = 0.1027738338969461
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Tnfsf13"
$ws.Cells.Item(8, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1501973333333333
$ws.Cells.Item(8, 8).Value = 0.450592
$ws.Cells.Item(8, 9).Value = 0.03220038283833477
$ws.Cells.Item(8, 10).Value = 0.03220038283833477
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.2347256666666667
$ws.Cells.Item(8, 14).Value = 0.7041770000000001
$ws.Cells.Item(8, 15).Value = 0.02963901389354162
$ws.Cells.Item(8, 16).Value = 0.02963901389354162
$ws.Cells.Item(8, 17).Value = 0.03525516919822223
$ws.Cells.Item(8, 18).Value = 0.317296522784
$ws.Cells.Item(8, 19).Value = 0.0009543875943227635
$ws.Cells.Item(8, 20).Value = 0.0009543875943227635
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Tnfsf13"
$ws.Cells.Item(9, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1501973333333333
$ws.Cells.Item(9, 8).Value = 0.450592
$ws.Cells.Item(9, 9).Value = 0.03220038283833477
$ws.Cells.Item(9, 10).Value = 0.03220038283833477
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.035534
$ws.Cells.Item(9, 14).Value = 0.106602
$ws.Cells.Item(9, 15).Value = 0.004486909057068498
$ws.Cells.Item(9, 16).Value = 0.004486909057068498
$ws.Cells.Item(9, 17).Value = 0.005337112042666667
$ws.Cells.Item(9, 18).Value = 0.048034008384
$ws.Cells.Item(9, 19).Value = 0.0001444801893983973
$ws.Cells.Item(9, 20).Value = 0.0001444801893983973
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Tnfsf13"
$ws.Cells.Item(10, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1501973333333333
$ws.Cells.Item(10, 8).Value = 0.450592
$ws.Cells.Item(10, 9).Value = 0.03220038283833477
$ws.Cells.Item(10, 10).Value = 0.03220038283833477
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.649223333333333
$ws.Cells.Item(10, 14).Value = 22.94767
$ws.Cells.Item(10, 15).Value = 0.9658740770493899
$ws.Cells.Item(10, 16).Value = 0.9658740770493899
$ws.Cells.Item(10, 17).Value = 1.148892946737778
$ws.Cells.Item(10, 18).Value = 10.34003652064
$ws.Cells.Item(10, 19).Value = 0.03110151505461361
$ws.Cells.Item(10, 20).Value = 0.03110151505461361
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Tnfsf13"
$ws.Cells.Item(11, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.675030666666667
$ws.Cells.Item(11, 8).Value = 8.025092000000001
$ws.Cells.Item(11, 9).Value = 0.573492282847582
$ws.Cells.Item(11, 10).Value = 0.573492282847582
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.2347256666666667
$ws.Cells.Item(11, 14).Value = 0.7041770000000001
$ws.Cells.Item(11, 15).Value = 0.02963901389354162
$ws.Cells.Item(11, 16).Value = 0.02963901389354162
$ws.Cells.Item(11, 17).Value = 0.6278983565871112
$ws.Cells.Item(11, 18).Value = 5.651085209284001
$ws.Cells.Item(11, 19).Value = 0.01699774573915839
$ws.Cells.Item(11, 20).Value = 0.01699774573915839
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Tnfsf13"
$ws.Cells.Item(12, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.675030666666667
$ws.Cells.Item(12, 8).Value = 8.025092000000001
$ws.Cells.Item(12, 9).Value = 0.573492282847582
$ws.Cells.Item(12, 10).Value = 0.573492282847582
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.035534
$ws.Cells.Item(12, 14).Value = 0.106602
$ws.Cells.Item(12, 15).Value = 0.004486909057068498
$ws.Cells.Item(12, 16).Value = 0.004486909057068498
$ws.Cells.Item(12, 17).Value = 0.09505453970933335
$ws.Cells.Item(12, 18).Value = 0.8554908573840001
$ws.Cells.Item(12, 19).Value = 0.002573207718067705
$ws.Cells.Item(12, 20).Value = 0.002573207718067705
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Tnfsf13"
$ws.Cells.Item(13, 3).Value = "Tnfrsf13b"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.675030666666667
$ws.Cells.Item(13, 8).Value = 8.025092000000001
$ws.Cells.Item(13, 9).Value = 0.573492282847582
$ws.Cells.Item(13, 10).Value = 0.573492282847582
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 7.649223333333333
$ws.Cells.Item(13, 14).Value = 22.94767
$ws.Cells.Item(13, 15).Value = 0.9658740770493899
$ws.Cells.Item(13, 16).Value = 0.9658740770493899
$ws.Cells.Item(13, 17).Value = 20.46190699284889
$ws.Cells.Item(13, 18).Value = 184.15716293564
$ws.Cells.Item(13, 19).Value = 0.553921329390356
$ws.Cells.Item(13, 20).Value = 0.553921329390356
